$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 614, pushing existing rows 614+ down to 616+
$ws.Rows("614:615").Insert()

# --- New row 614 ---
$ws.Cells.Item(614,1).Value = 3
$ws.Cells.Item(614,2).Value = "Femacal de La Calera"
$ws.Cells.Item(614,3).Value = "Coquimbo"
$ws.Cells.Item(614,4).Value = 44706
$ws.Cells.Item(614,5).Value = 5
$ws.Cells.Item(614,6).Value = "Fruta"
$ws.Cells.Item(614,7).Value = 100108
$ws.Cells.Item(614,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(614,9).Value = 100108005
$ws.Cells.Item(614,10).Value = "Piña"
$ws.Cells.Item(614,11).Value = "Caramelo"
$ws.Cells.Item(614,12).Value = "Primera"
$ws.Cells.Item(614,13).Value = 54
$ws.Cells.Item(614,14).Value = 18000
$ws.Cells.Item(614,15).Value = 18000
$ws.Cells.Item(614,16).Value = 18000
$ws.Cells.Item(614,17).Value = "$/caja 12 unidades"
$ws.Cells.Item(614,18).Value = "Ecuador"
$ws.Cells.Item(614,19).Value = 1500
$ws.Cells.Item(614,20).Value = 12

# --- New row 615 ---
$ws.Cells.Item(615,1).Value = 3
$ws.Cells.Item(615,2).Value = "Femacal de La Calera"
$ws.Cells.Item(615,3).Value = "Coquimbo"
$ws.Cells.Item(615,4).Value = 44706
$ws.Cells.Item(615,5).Value = 5
$ws.Cells.Item(615,6).Value = "Fruta"
$ws.Cells.Item(615,7).Value = 100108
$ws.Cells.Item(615,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(615,9).Value = 100108005
$ws.Cells.Item(615,10).Value = "Piña"
$ws.Cells.Item(615,11).Value = "Caramelo"
$ws.Cells.Item(615,12).Value = "Segunda"
$ws.Cells.Item(615,13).Value = 108
$ws.Cells.Item(615,14).Value = 18000
$ws.Cells.Item(615,15).Value = 18000
$ws.Cells.Item(615,16).Value = 18000
$ws.Cells.Item(615,17).Value = "$/caja 14 unidades"
$ws.Cells.Item(615,18).Value = "Ecuador"
$ws.Cells.Item(615,19).Value = 1286
$ws.Cells.Item(615,20).Value = 14
